# "update manufacturer structure to include notice field"
#
# The Manufacturer block (originally the merged header R1:S1 with two data
# columns R:S) gains a third data column for a new "notice" field. In
# Excel terms this is a plain column insert immediately before the old
# "Cost:" block (column T), which:
#   - shifts the old Cost:/Icon:/Description: columns (T:V, W, X) one
#     column to the right (now U:W, X, Y) together with their formatting,
#   - pushes the dimension/used-range out to column Y,
#   - and leaves the brand-new column (T) carrying a copy of the
#     formatting of the column immediately to its left (S), i.e. the
#     Manufacturer block's look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new column before "T" (the old "Cost:" header started there).
$ws.Columns("T:T").Insert()

# Re-merge the "Manufacturer:" header so it spans its new third column
# (was R1:S1, now R1:T1), centered like the sheet's other group headers.
$ws.Range("R1:S1").UnMerge()
$ws.Range("R1:T1").Merge()
$ws.Range("R1:T1").HorizontalAlignment = -4108

# Populate the new "notice" field in row 2 under Manufacturer.
$ws.Range("T2").Value = "notice"

# Restore the author's selection in the saved view.
$ws.Range("U7").Select()
